# Override the "Transporte" class for a handful of point_index groups now
# that the analysis buffer = 75 m. These rows were mis-classified as
# "Transporte" and should instead read "Outros" / "Urbanizado" in column C
# (Uso_solo_simplificado), matching the reclassification performed upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# point_index group 18 (rows 74-77): Transporte -> Outros
$ws.Range("C74:C77").Value = "Outros"

# point_index groups 27, 36, 37, 43, 44: Transporte -> Urbanizado
$ws.Range("C110:C113").Value = "Urbanizado"
$ws.Range("C146:C153").Value = "Urbanizado"
$ws.Range("C174:C181").Value = "Urbanizado"
